# Regenerate save_data to use K (strikeout count) instead of Strike#.
# Update column G ("K") values for rows 2-9 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 2
    6 = 1
    7 = 0
    8 = 0
    9 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
